$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F13").Value = 76
$ws.Range("F14").Value = 52
$ws.Range("F15").Value = 613
$ws.Range("F17").Value = 7161
$ws.Range("F19").Value = 7473
$ws.Range("F21").Value = 56993
$ws.Range("F22").Value = 4608
$ws.Range("F23").Value = 1046
$ws.Range("F24").Value = 901
$ws.Range("F25").Value = 463
$ws.Range("F30").Value = 4554
$ws.Range("F31").Value = 591
$ws.Range("F32").Value = 78
$ws.Range("F33").Value = 36
$ws.Range("F35").Value = 1277
$ws.Range("F36").Value = 1400
$ws.Range("F44").Value = 213
$ws.Range("F46").Value = 183

$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value = 15
$ws.Range("F10").Value = 45
$ws.Range("F11").Value = 7536
$ws.Range("F24").Value = 22
$ws.Range("F36").Value = 45
$ws.Range("F40").Value = 6
$ws.Range("F41").Value = 99
$ws.Range("F48").Value = 152

$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value = 2327
$ws.Range("F9").Value = 9366
$ws.Range("F10").Value = 1713
$ws.Range("F11").Value = 168
$ws.Range("F12").Value = 97
$ws.Range("F15").Value = 228
$ws.Range("F16").Value = 2070
$ws.Range("F17").Value = 426

$ws = $wb.Worksheets.Item(4)
$ws.Range("F7").Value = 1713
$ws.Range("F8").Value = 168
$ws.Range("F9").Value = 97
$ws.Range("F13").Value = 228
$ws.Range("F14").Value = 2070
$ws.Range("F15").Value = 52
$ws.Range("F16").Value = 613
$ws.Range("F18").Value = 7161
$ws.Range("F19").Value = 56993
$ws.Range("F22").Value = 4608
$ws.Range("F23").Value = 1046
$ws.Range("F24").Value = 463
$ws.Range("F27").Value = 4554
$ws.Range("F28").Value = 591
$ws.Range("F29").Value = 79
$ws.Range("F31").Value = 1277
$ws.Range("F33").Value = 427
$ws.Range("F41").Value = 22
$ws.Range("F44").Value = 183
$ws.Range("F47").Value = 45
